# RoomOccupancy / Wydziały.xlsx
# The faculty "Wydział Zastosowań Informatyki I Matematyki" had its short
# name (abbreviation) changed from "WZIM" to "WZIiM" (row 14, column B of
# the lookup table on the single worksheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = "WZIiM"

# Reflect the author's last cursor position when the file was saved.
$ws.Range("B18").Select()
